# Actualización desde MV -datos-
# Adds four new daily rows (02-11-2021 .. 05-11-2021) to the bottom of the
# "Swap promedio camara" data sheet, mirroring the existing row layout
# (columns A..M): A = date label (text), B..M = numeric rates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data exactly as it appears in the target workbook.
$rows = @(
    @{ Row = 213; Date = "02-11-2021"; Values = @(5.21, 5.25, 5.33, 5.42, 5.64, -0.8100000000000001, 0.2, 0.9, 1.23, 1.46, 1.95, 2.01) },
    @{ Row = 214; Date = "03-11-2021"; Values = @(5.22, 5.28, 5.38, 5.47, 5.69, -0.88, 0.18, 0.89, 1.23, 1.47, 1.91, 1.97) },
    @{ Row = 215; Date = "04-11-2021"; Values = @(5.23, 5.28, 5.38, 5.48, 5.73, -1, 0.13, 0.85, 1.2, 1.43, 1.85, 1.91) },
    @{ Row = 216; Date = "05-11-2021"; Values = @(5.14, 5.18, 5.29, 5.4, 5.64, -1.02, 0.11, 0.83, 1.17, 1.39, 1.8, 1.87) }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    # Column A holds a date-shaped label that must stay plain text (just
    # like the other ~212 rows above it). A direct .Value assignment would
    # be auto-recognised as a real date, so instead write it as a literal
    # text formula and collapse it back down to a plain cached value via
    # copy / paste-values (keeps the cell unstyled, same as its neighbours).
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '="' + $entry.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    # Columns B..M are plain numeric values.
    $col = 2
    foreach ($val in $entry.Values) {
        $ws.Cells.Item($r, $col).Value = $val
        $col = $col + 1
    }
}

$excel.CutCopyMode = 0
